$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the entry for Thursday 23.10. (row 38): worked 7 hours, 17:00 - 24:00
$ws.Range("D38").Value = 7
$ws.Range("E38").Value = "17:00 - 24:00"

# Restore the view scroll position / selection as saved
$ws.Range("I43").Select()
$ws.Application.ActiveWindow.ScrollRow = 29
